# Apply the content-rotation edit to LOB1276.docx
$d = $word.ActiveDocument

# ---- Phase 1: replace each source block (in its current slot) with a unique placeholder token ----
# (Placeholder tokens are needed because the edit is a cyclic rotation of text blocks
#  through fixed paragraph slots; a direct sequential replace would clobber data.)
$search = "Propiciar o entendimento sobre os conceitos básicos, teoria, metodologias de análise e práticas sobre a poluição do solo e águas subterrâneas."
$null = $d.Content.Find.Execute($search, $true, $false, $false, $false, $false, $true, 1, $false, "@@TOKEN0@@", 2)

$search = "5464150 - Mariana Consiglio Kasemodel"
$null = $d.Content.Find.Execute($search, $true, $false, $false, $false, $false, $true, 1, $false, "@@TOKEN1@@", 2)

$search = "BOSCOV, M.E.G. Geotecnia Ambiental. Oficina de Textos, 2008. 248 p." + [char]11 + "CETESB. Decisão da Diretoria N. 38/2017/C, de 7 de fevereiro de 2017. 65 p." + [char]11 + "OLIVEIRA, A. M. S.; JERÔNIMO, J. Geologia de Engenharia e Ambiental, ABGE, 2018. 912 p." + [char]11 + "SHARMA, H. D.; REDDY, K. R. Geoenvironmental engineering, Wiley, 2004. 992p. " + [char]11 + "YONG, R. N. Geoenvironmental engineering. contaminated soils, Pollutant fate and migration. CRC Press, 2001. 307p."
$null = $d.Content.Find.Execute($search, $true, $false, $false, $false, $false, $true, 1, $false, "@@TOKEN2@@", 2)

$search = "1 (uma) prova escrita"
$null = $d.Content.Find.Execute($search, $true, $false, $false, $false, $false, $true, 1, $false, "@@TOKEN3@@", 2)

$search = "Média ponderada de provas e atividades."
$null = $d.Content.Find.Execute($search, $true, $false, $false, $false, $false, $true, 1, $false, "@@TOKEN4@@", 2)

$search = "Aulas teóricas expositivas, atividades individuais e em grupo, relatórios e provas."
$null = $d.Content.Find.Execute($search, $true, $false, $false, $false, $false, $true, 1, $false, "@@TOKEN5@@", 2)

$search = "Conceitos introdutórios: composição química e mineralógica do solo, horizontes de solo, unidades aquíferas" + [char]11 + "Principais poluentes dos solos e águas subterrâneas" + [char]11 + "Legislação: valores orientadores" + [char]11 + "Interação entre solo e poluentes " + [char]11 + "Fluxo de água em zonas não saturadas e saturadas" + [char]11 + "Fluxo de água em fraturas de rocha" + [char]11 + "Transporte de substâncias miscíveis na água no solo e águas subterrâneas: mecanismos, modelos e soluções " + [char]11 + "Transporte de substâncias não miscíveis na água (NAPL) no solo e águas subterrâneas: mecanismos, modelos e soluções " + [char]11 + "Gerenciamento de áreas contaminadas: investigação preliminar, técnicas de investigação geológica-geotécnica; modelo conceitual; análise de risco e técnicas de intervenção " + [char]11 + "Exemplos de aplicação em problemas geoambientais" + [char]11 + "A disciplina pode contar com viagens didáticas para complementação do conteúdo da disciplina"
$null = $d.Content.Find.Execute($search, $true, $false, $false, $false, $false, $true, 1, $false, "@@TOKEN6@@", 2)

$search = "Principais fatores condicionantes atuantes em processos de contaminação de solos e águas subterrânea. Técnicas de investigação, monitoramento, contenção e recuperação em geotecnia ambiental."
$null = $d.Content.Find.Execute($search, $true, $false, $false, $false, $false, $true, 1, $false, "@@TOKEN7@@", 2)

# ---- Phase 2: replace each placeholder token with its final destination content ----
$replacement = "Principais fatores condicionantes atuantes em processos de contaminação de solos e águas subterrânea. Técnicas de investigação, monitoramento, contenção e recuperação em geotecnia ambiental."
$null = $d.Content.Find.Execute("@@TOKEN0@@", $true, $false, $false, $false, $false, $true, 1, $false, $replacement, 2)

$replacement = "Propiciar o entendimento sobre os conceitos básicos, teoria, metodologias de análise e práticas sobre a poluição do solo e águas subterrâneas."
$null = $d.Content.Find.Execute("@@TOKEN1@@", $true, $false, $false, $false, $false, $true, 1, $false, $replacement, 2)

$replacement = "5464150 - Mariana Consiglio Kasemodel"
$null = $d.Content.Find.Execute("@@TOKEN2@@", $true, $false, $false, $false, $false, $true, 1, $false, $replacement, 2)

$replacement = "BOSCOV, M.E.G. Geotecnia Ambiental. Oficina de Textos, 2008. 248 p." + [char]11 + "CETESB. Decisão da Diretoria N. 38/2017/C, de 7 de fevereiro de 2017. 65 p." + [char]11 + "OLIVEIRA, A. M. S.; JERÔNIMO, J. Geologia de Engenharia e Ambiental, ABGE, 2018. 912 p." + [char]11 + "SHARMA, H. D.; REDDY, K. R. Geoenvironmental engineering, Wiley, 2004. 992p. " + [char]11 + "YONG, R. N. Geoenvironmental engineering. contaminated soils, Pollutant fate and migration. CRC Press, 2001. 307p."
$null = $d.Content.Find.Execute("@@TOKEN3@@", $true, $false, $false, $false, $false, $true, 1, $false, $replacement, 2)

$replacement = "1 (uma) prova escrita"
$null = $d.Content.Find.Execute("@@TOKEN4@@", $true, $false, $false, $false, $false, $true, 1, $false, $replacement, 2)

$replacement = "Média ponderada de provas e atividades."
$null = $d.Content.Find.Execute("@@TOKEN5@@", $true, $false, $false, $false, $false, $true, 1, $false, $replacement, 2)

$replacement = "Aulas teóricas expositivas, atividades individuais e em grupo, relatórios e provas."
$null = $d.Content.Find.Execute("@@TOKEN6@@", $true, $false, $false, $false, $false, $true, 1, $false, $replacement, 2)

$replacement = "Conceitos introdutórios: composição química e mineralógica do solo, horizontes de solo, unidades aquíferas" + [char]11 + "Principais poluentes dos solos e águas subterrâneas" + [char]11 + "Legislação: valores orientadores" + [char]11 + "Interação entre solo e poluentes " + [char]11 + "Fluxo de água em zonas não saturadas e saturadas" + [char]11 + "Fluxo de água em fraturas de rocha" + [char]11 + "Transporte de substâncias miscíveis na água no solo e águas subterrâneas: mecanismos, modelos e soluções " + [char]11 + "Transporte de substâncias não miscíveis na água (NAPL) no solo e águas subterrâneas: mecanismos, modelos e soluções " + [char]11 + "Gerenciamento de áreas contaminadas: investigação preliminar, técnicas de investigação geológica-geotécnica; modelo conceitual; análise de risco e técnicas de intervenção " + [char]11 + "Exemplos de aplicação em problemas geoambientais" + [char]11 + "A disciplina pode contar com viagens didáticas para complementação do conteúdo da disciplina"
$null = $d.Content.Find.Execute("@@TOKEN7@@", $true, $false, $false, $false, $false, $true, 1, $false, $replacement, 2)

# ---- OBJ_EN <-> RESUMO_EN swap (independent 2-cycle) ----
$search = "Provide knowledge of the basic concepts, theory, analysis methodologies of soil and groundwater pollution."
$null = $d.Content.Find.Execute($search, $true, $false, $false, $false, $false, $true, 1, $false, "@@TOKEN_EN_A@@", 2)

$search = "Main aspects acting in processes such as erosion, gravitational mass movements and contamination of geological materials and groundwater. Investigation, monitoring, containment and recuperation techniques in environmental geotechnics"
$null = $d.Content.Find.Execute($search, $true, $false, $false, $false, $false, $true, 1, $false, "@@TOKEN_EN_B@@", 2)

$replacement = "Main aspects acting in processes such as erosion, gravitational mass movements and contamination of geological materials and groundwater. Investigation, monitoring, containment and recuperation techniques in environmental geotechnics"
$null = $d.Content.Find.Execute("@@TOKEN_EN_A@@", $true, $false, $false, $false, $false, $true, 1, $false, $replacement, 2)

$replacement = "Provide knowledge of the basic concepts, theory, analysis methodologies of soil and groundwater pollution."
$null = $d.Content.Find.Execute("@@TOKEN_EN_B@@", $true, $false, $false, $false, $false, $true, 1, $false, $replacement, 2)

Write-Output "done"